$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Embalagem Panilha")

$rows = @(
    @("ENTRADA","dsa","das","das","das","EMBALAGEM","das","das","das","das","das","das","das","das"),
    @("ENTRADA","asdasx","saxsa","xsa","xas","EMBALAGEM","xsa","xas","xsa","xas","xas","xas","xas","xas"),
    @("ENTRADA","asdasx","saxsa","xsa","xas","EMBALAGEM","xsa","xas","xsa","xas","xas","xas","xas","xas"),
    @("ENTRADA","df","df","df","qwe","EMBALAGEM","df","qwe","qwe","qwe","qwe","qwe","qwe","qwe"),
    @("ENTRADA","df","df","df","asd","EMBALAGEM","df","asd","asd","asd","asd","asd","asd","asd"),
    @("ENTRADA","df","df","df","rfv","EMBALAGEM","df","rfv","rfv","rfv","frv","rfv","rfv","rfv"),
    @("ENTRADA","df","df","df","yhb","EMBALAGEM","df","yhb","yhb","yhb","yhb","yhb","yhb","yhb"),
    @("ENTRADA","df","df","df","ikm","EMBALAGEM","df","ikmik","ikm","ikm","ikm","ikm","ik","ikm"),
    @("ENTRADA","df","df","df","dfg","EMBALAGEM","df","dfg","dfg","dfg","dfg","dfg","dfg","dfg")
)

$startRow = 3
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowData = $rows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
